$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full 3x3 grid of sending/target clusters with updated metrics.
$data = @(
    @{ Row=2;  A="ECs";  D="ECs";  E=3; F=1; G=203.7816646666667; H=611.344994;        I=0.6667327591988204; J=0.6667327591988205; K=3; L=1;                  M=0.686889;           N=2.060667;  O=0.2801468581979037; P=0.2801468581979037; Q=139.975383861222;  R=1259.778454750998;  S=0.186783087747169;   T=0.1867830877471691 },
    @{ Row=3;  A="ECs";  D="FAPs"; E=3; F=1; G=203.7816646666667; H=611.344994;        I=0.6667327591988204; J=0.6667327591988205; K=2; L=0.6666666666666666; M=0.7012299999999999; N=2.10369;   O=0.2859958179183478; P=0.2859958179183478; Q=142.8978167142066; R=1286.08035042786;   S=0.1906827808000235;  T=0.1906827808000235 },
    @{ Row=4;  A="ECs";  D="sCs";  E=3; F=1; G=203.7816646666667; H=611.344994;        I=0.6667327591988204; J=0.6667327591988205; K=3; L=1;                  M=1.06377;            N=3.19131;   O=0.4338573238837485; P=0.4338573238837484; Q=216.77682142246;   R=1950.99139280214;   S=0.2892668906516279;  T=0.2892668906516279 },
    @{ Row=5;  A="FAPs"; D="ECs";  E=3; F=1; G=63.14058933333333; H=189.421768;        I=0.2065833519051582; J=0.2065833519051582; K=3; L=1;                  M=0.686889;           N=2.060667;  O=0.2801468581979037; P=0.2801468581979037; Q=43.370576266584;   R=390.335186399256;   S=0.057873676992222;   T=0.057873676992222 },
    @{ Row=6;  A="FAPs"; D="FAPs"; E=3; F=1; G=63.14058933333333; H=189.421768;        I=0.2065833519051582; J=0.2065833519051582; K=2; L=0.6666666666666666; M=0.7012299999999999; N=2.10369;   O=0.2859958179183478; P=0.2859958179183478; Q=44.27607545821333; R=398.4846791239199;  S=0.05908197469642959; T=0.0590819746964296 },
    @{ Row=7;  A="FAPs"; D="sCs";  E=3; F=1; G=63.14058933333333; H=189.421768;        I=0.2065833519051582; J=0.2065833519051582; K=3; L=1;                  M=1.06377;            N=3.19131;   O=0.4338573238837485; P=0.4338573238837484; Q=67.16706471512001; R=604.50358243608;    S=0.0896277002165066;  T=0.08962770021650658 },
    @{ Row=8;  A="sCs";  D="ECs";  E=3; F=1; G=38.719942;         H=116.159826;        I=0.1266838888960214; J=0.1266838888960214; K=3; L=1;                  M=0.686889;           N=2.060667;  O=0.2801468581979037; P=0.2801468581979037; Q=26.596302240438;   R=239.366720163942;   S=0.0354900934585127;  T=0.03549009345851271 },
    @{ Row=9;  A="sCs";  D="FAPs"; E=3; F=1; G=38.719942;         H=116.159826;        I=0.1266838888960214; J=0.1266838888960214; K=2; L=0.6666666666666666; M=0.7012299999999999; N=2.10369;   O=0.2859958179183478; P=0.2859958179183478; Q=27.15158492865999; R=244.36426435794;    S=0.03623106242189474; T=0.03623106242189475 },
    @{ Row=10; A="sCs";  D="sCs";  E=3; F=1; G=38.719942;         H=116.159826;        I=0.1266838888960214; J=0.1266838888960214; K=3; L=1;                  M=1.06377;            N=3.19131;   O=0.4338573238837485; P=0.4338573238837484; Q=41.18911270134;    R=370.70201431206;    S=0.05496273301561396; T=0.05496273301561396 }
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = "Gnai2"
    $ws.Cells.Item($r, 3).Value = "Tbxa2r"
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 9).Value = $d.I
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = $d.Q
    $ws.Cells.Item($r, 18).Value = $d.R
    $ws.Cells.Item($r, 19).Value = $d.S
    $ws.Cells.Item($r, 20).Value = $d.T
}
